$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '28.435.27'
$ws.Range("E2").Value = '  +0.01%  '
Set-TextValue "D3" '1.551.54'
$ws.Range("E3").Value = '  -2.04%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.24%  '
Set-TextValue "D5" '210.59'
$ws.Range("E5").Value = '  -1.41%  '
$ws.Range("E6").Value = '  -1.61%  '
Set-TextValue "D7" '1.00'
$ws.Range("E7").Value = '  -0.26%  '
Set-TextValue "D8" '24.01'
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("E10").Value = '  -1.42%  '
Set-TextValue "D11" '0.0890'
$ws.Range("E11").Value = '  -0.44%  '
Set-TextValue "D12" '1.772.80'
$ws.Range("E12").Value = '  -2.09%  '
Set-TextValue "D13" '1.550.81'
$ws.Range("E13").Value = '  -2.14%  '
Set-TextValue "D14" '28.430.85'
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  -2.16%  '
Set-TextValue "D16" '0.510'
$ws.Range("E16").Value = '  -1.98%  '
Set-TextValue "D17" '60.98'
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("E19").Value = '  -1.76%  '
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("E24").Value = '  -2.05%  '
Set-TextValue "D25" '151.20'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("E26").Value = '  -1.82%  '
$ws.Range("E27").Value = '  -1.32%  '
Set-TextValue "D28" '1.00'
Set-TextValue "D29" '6.22'
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("E30").Value = '  -3.32%  '
$ws.Range("E31").Value = '  -4.39%  '
$ws.Range("E32").Value = '  -2.15%  '
Set-TextValue "D33" '1.380.57'
$ws.Range("E33").Value = '  -1.27%  '
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("E36").Value = '  -3.17%  '
Set-TextValue "D37" '2.30'
$ws.Range("E37").Value = '  -2.89%  '
Set-TextValue "D38" '2.58'
$ws.Range("E38").Value = '  -3.12%  '
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D40" '0.511'
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D41" '1.91'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("E43").Value = '  -2.42%  '
$ws.Range("E44").Value = '  -0.87%  '
Set-TextValue "D45" '5.37'
$ws.Range("E45").Value = '  -1.19%  '
Set-TextValue "D46" '61.80'
$ws.Range("E46").Value = '  -1.99%  '
Set-TextValue "D47" '1.685.61'
$ws.Range("E47").Value = '  -2.11%  '
Set-TextValue "D48" '0.875'
$ws.Range("E48").Value = '  -9.09%  '
Set-TextValue "D49" '85.18'
$ws.Range("E49").Value = '  -1.70%  '
Set-TextValue "D50" '42.99'
$ws.Range("E50").Value = '  +8.40%  '
$ws.Range("E51").Value = '  -2.55%  '
